# BARD1 Filtering bug fix
# The start and end coordinates in the "gene" sheet input file were
# swapped; fix by swapping the A2 (start) and B2 (end) values back.

$wb = $excel.ActiveWorkbook

$geneSheet = $wb.Worksheets.Item("gene")

$startVal = $geneSheet.Range("A2").Value2
$endVal = $geneSheet.Range("B2").Value2

$geneSheet.Range("A2").Value2 = $endVal
$geneSheet.Range("B2").Value2 = $startVal

# Reflect the user's final navigation state: "gene" sheet active,
# selection on C9.
$geneSheet.Activate()
$geneSheet.Range("C9").Select()
